$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.309.50"
$ws.Range("E2").Value = "'  +0.33%  "
$ws.Range("D3").Value = "'1.868.53"
$ws.Range("E3").Value = "'  +0.44%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'235.13"
$ws.Range("E5").Value = "'  -0.36%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.14%  "
$ws.Range("D7").Value = "'0.4701"
$ws.Range("E7").Value = "'  +0.61%  "
$ws.Range("D8").Value = "'0.2867"
$ws.Range("E8").Value = "'  +0.43%  "
$ws.Range("D9").Value = "'0.06575"
$ws.Range("E9").Value = "'  +0.84%  "
$ws.Range("D10").Value = "'21.67"
$ws.Range("E10").Value = "'  -0.41%  "
$ws.Range("D11").Value = "'0.08027"
$ws.Range("E11").Value = "'  +1.70%  "
$ws.Range("D12").Value = "'97.01"
$ws.Range("E12").Value = "'  -0.33%  "
$ws.Range("D13").Value = "'1.872.00"
$ws.Range("E13").Value = "'  +0.30%  "
$ws.Range("D14").Value = "'5.115"
$ws.Range("E14").Value = "'  -0.91%  "
$ws.Range("D15").Value = "'0.6848"
$ws.Range("E15").Value = "'  +0.85%  "
$ws.Range("D16").Value = "'268.77"
$ws.Range("D17").Value = "'30.329.71"
$ws.Range("E17").Value = "'  +0.42%  "
$ws.Range("D18").Value = "'13.95"
$ws.Range("E18").Value = "'  +3.45%  "
$ws.Range("D19").Value = "'0.000007645"
$ws.Range("E19").Value = "'  +4.73%  "
$ws.Range("E20").Value = "'  +0.13%  "
$ws.Range("D21").Value = "'2.116.76"
$ws.Range("E21").Value = "'  +0.37%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.21%  "
$ws.Range("D23").Value = "'5.270"
$ws.Range("E23").Value = "'  -1.81%  "
$ws.Range("D24").Value = "'6.206"
$ws.Range("E24").Value = "'  +0.72%  "
$ws.Range("D25").Value = "'9.406"
$ws.Range("E25").Value = "'  +2.06%  "
$ws.Range("D26").Value = "'168.73"
$ws.Range("E26").Value = "'  +0.60%  "
$ws.Range("D27").Value = "'18.87"
$ws.Range("E27").Value = "'  -0.93%  "
$ws.Range("D28").Value = "'1.949"
$ws.Range("E28").Value = "'  +1.10%  "
$ws.Range("D29").Value = "'1.371"
$ws.Range("E29").Value = "'  -0.64%  "
$ws.Range("D30").Value = "'0.09870"
$ws.Range("E30").Value = "'  +1.55%  "
$ws.Range("D31").Value = "'4.372"
$ws.Range("E31").Value = "'  +0.25%  "
$ws.Range("E32").Value = "'  -1.20%  "
$ws.Range("D33").Value = "'4.067"
$ws.Range("E33").Value = "'  +0.55%  "
$ws.Range("D34").Value = "'0.04694"
$ws.Range("E34").Value = "'  -0.63%  "
$ws.Range("D35").Value = "'1.134"
$ws.Range("E35").Value = "'  +0.12%  "
$ws.Range("D36").Value = "'0.7006"
$ws.Range("E36").Value = "'  -1.00%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("E37").Value = "'  -0.14%  "
$ws.Range("D38").Value = "'0.01871"
$ws.Range("E38").Value = "'  +0.49%  "
$ws.Range("D39").Value = "'2.628"
$ws.Range("E39").Value = "'  -0.11%  "
$ws.Range("D40").Value = "'6.279"
$ws.Range("E40").Value = "'  -0.70%  "
$ws.Range("D41").Value = "'72.18"
$ws.Range("E41").Value = "'  -2.99%  "
$ws.Range("D42").Value = "'1.957"
$ws.Range("E42").Value = "'  +0.42%  "
$ws.Range("D43").Value = "'0.8419"
$ws.Range("E43").Value = "'  -0.68%  "
$ws.Range("D44").Value = "'0.4166"
$ws.Range("E44").Value = "'  -0.12%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "'  +0.10%  "
$ws.Range("D46").Value = "'102.92"
$ws.Range("E46").Value = "'  -0.37%  "
$ws.Range("D47").Value = "'9.166"
$ws.Range("E47").Value = "'  -0.94%  "
$ws.Range("D48").Value = "'7.046"
$ws.Range("E48").Value = "'  -1.80%  "
$ws.Range("D49").Value = "'912.90"
$ws.Range("E49").Value = "'  -5.87%  "
$ws.Range("D50").Value = "'34.52"
$ws.Range("D51").Value = "'0.05684"
$ws.Range("E51").Value = "'  +0.85%  "
